$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.277.36'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.664.79'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.93'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5309'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.010'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2639'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +1.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06362'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.57'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +0.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07859'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.657.99'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.893.09'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5528'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8174'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.65'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.341.69'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.010'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.677'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +2.51%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.52'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.23'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.039'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.012'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.44'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +1.74%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1227'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -1.77%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.215'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.15'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.482'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +3.29%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05945'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.282'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.591'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +1.96%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.285'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.614'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +2.64%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9620'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.828'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +1.05%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.425'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5812'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +2.76%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01606'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8660'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +2.07%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.870'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +1.02%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.048.16'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +2.33%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '104.12'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +1.50%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.805.91'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.46'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈106'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -5.65%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.015'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4381'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +2.22%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.975'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +2.00%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05163'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +0.23%  '
